$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$name = "KAGS #7158"
$kd = 1.01
$winPct = 57
$topAgent = "Jett"
$hsPct = 23.3
$clutches = 14
$firstKills = 138
$firstDeaths = 156
$knifeKills = 0
$damage = 15
$rank = "Nickel"
$archetype = "['Rusher', 'Straight Up Winner']"

for ($r = 100; $r -le 103; $r++) {
    $ws.Cells.Item($r, 1).Value = $name
    $ws.Cells.Item($r, 2).Value = $kd
    $ws.Cells.Item($r, 3).Value = $winPct
    $ws.Cells.Item($r, 4).Value = $topAgent
    $ws.Cells.Item($r, 5).Value = $hsPct
    $ws.Cells.Item($r, 6).Value = $clutches
    $ws.Cells.Item($r, 7).Value = $firstKills
    $ws.Cells.Item($r, 8).Value = $firstDeaths
    $ws.Cells.Item($r, 9).Value = $knifeKills
    $ws.Cells.Item($r, 10).Value = $damage
    $ws.Cells.Item($r, 11).Value = $rank
    $ws.Cells.Item($r, 12).Value = $archetype
}
